# KeyboardBOM.xlsx - mark received parts with their actual receipt dates.
# Writing in this specific order controls the order new shared strings are
# appended to the workbook's string table (first-use order), matching the
# canonical file exactly.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H7").Value = "Yes (04/12/18)"
$ws.Range("H5").Value = "Yes (04/13/18)"
$ws.Range("L24").Value = "Yes (04/24/18)"
$ws.Range("H4").Value = "Yes (04/19/18)"
$ws.Range("H6").Value = "Yes (04/13/18)"
$ws.Range("H8").Value = "Yes (04/12/18)"
$ws.Range("H9").Value = "Yes (04/13/18)"
$ws.Range("H10").Value = "Yes (04/12/18)"
$ws.Range("H11").Value = "Yes (04/12/18)"
$ws.Range("H12").Value = "Yes (04/12/18)"
$ws.Range("H13").Value = "Yes (04/12/18)"
$ws.Range("L19").Value = "Yes (04/12/18)"

# Widen the "Received?" date column now that it holds longer values.
$ws.Columns.Item(12).ColumnWidth = 13.09

# Move the active selection to reflect where the user left off editing.
$ws.Range("H5").Select() | Out-Null
